$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing last-row value (idj2jd92j -> flirtest4)
$ws.Range("A2").Value = "flirtest4"

# Add new negative test rows (flirtest5 .. flirtest9), each with
# Automated / Test in columns B and C like the existing rows.
$names = @("flirtest5", "flirtest6", "flirtest7", "flirtest8", "flirtest9")
$row = 3
foreach ($name in $names) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = "Automated"
    $ws.Cells.Item($row, 3).Value = "Test"
    $row++
}

$ws.Range("B7").Select() | Out-Null
